$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet originally has data rows 2..223. New catastral-parcel records
# were discovered, so 16 new rows are appended (sheet grows to row 239) and
# a block of existing rows (158..223) is rewritten with the refreshed data.
# ---------------------------------------------------------------------------

# 1) Insert 16 blank rows at the bottom of the table (224..239), copying the
#    formatting (borders/bold/alignment on column A, etc.) from the last
#    existing data row (223) so the new rows look identical to the old ones.
$ws.Range("A223:C223").Copy()
$ws.Range("A224:C239").PasteSpecial(-4122)

# 2) Make sure column B keeps storing plain text (codice_particella values
#    such as ".507/2" or "7294" must not be re-interpreted as numbers).
$ws.Range("B158:B239").NumberFormat = "@"

$data = @(
  @{Row=158; A=156; B='.507/2'; C=413}
  @{Row=159; A=157; B='.507/3'; C=413}
  @{Row=160; A=158; B='.507/4'; C=413}
  @{Row=161; A=159; B='507/5'; C=413}
  @{Row=162; A=160; B='507/6'; C=413}
  @{Row=163; A=161; B='569/8'; C=413}
  @{Row=164; A=162; B='.569/5'; C=413}
  @{Row=165; A=163; B='7294'; C=413}
  @{Row=166; A=164; B='7293/2'; C=413}
  @{Row=167; A=165; B='7293/3'; C=413}
  @{Row=168; A=166; B='1.391'; C=413}
  @{Row=169; A=167; B='7292/1'; C=413}
  @{Row=170; A=168; B='1.358'; C=413}
  @{Row=171; A=169; B='1.359'; C=413}
  @{Row=172; A=170; B='7191/4'; C=413}
  @{Row=173; A=171; B='1.270'; C=413}
  @{Row=174; A=172; B='7618/3'; C=413}
  @{Row=175; A=173; B='7618/4'; C=413}
  @{Row=176; A=174; B='114'; C=413}
  @{Row=177; A=175; B='361'; C=413}
  @{Row=178; A=176; B='421'; C=413}
  @{Row=179; A=177; B='13094'; C=413}
  @{Row=180; A=178; B='7290/6'; C=413}
  @{Row=181; A=179; B='7288'; C=413}
  @{Row=182; A=180; B='45'; C=283}
  @{Row=183; A=181; B='84/1'; C=283}
  @{Row=184; A=182; B='137'; C=41}
  @{Row=185; A=183; B='140'; C=41}
  @{Row=186; A=184; B='272'; C=79}
  @{Row=187; A=185; B='1881/8'; C=79}
  @{Row=188; A=186; B='94/6'; C=251}
  @{Row=189; A=187; B='243'; C=253}
  @{Row=190; A=188; B='295'; C=253}
  @{Row=191; A=189; B='1403'; C=268}
  @{Row=192; A=190; B='316'; C=442}
  @{Row=193; A=191; B='53'; C=215}
  @{Row=194; A=192; B='454'; C=215}
  @{Row=195; A=193; B='420/80'; C=215}
  @{Row=196; A=194; B='420/92'; C=215}
  @{Row=197; A=195; B='420/93'; C=215}
  @{Row=198; A=196; B='420/94'; C=215}
  @{Row=199; A=197; B='420/95'; C=215}
  @{Row=200; A=198; B='420/96'; C=215}
  @{Row=201; A=199; B='420/97'; C=215}
  @{Row=202; A=200; B='420/101'; C=215}
  @{Row=203; A=201; B='420/102'; C=215}
  @{Row=204; A=202; B='420/106'; C=215}
  @{Row=205; A=203; B='420/107'; C=215}
  @{Row=206; A=204; B='420/109'; C=215}
  @{Row=207; A=205; B='420/110'; C=215}
  @{Row=208; A=206; B='705/11'; C=215}
  @{Row=209; A=207; B='756'; C=215}
  @{Row=210; A=208; B='798/3'; C=215}
  @{Row=211; A=209; B='1411/1'; C=256}
  @{Row=212; A=210; B='1411/2'; C=256}
  @{Row=213; A=211; B='1411/3'; C=256}
  @{Row=214; A=212; B='1411/4'; C=256}
  @{Row=215; A=213; B='1411/5'; C=256}
  @{Row=216; A=214; B='1412'; C=256}
  @{Row=217; A=215; B='1488'; C=256}
  @{Row=218; A=216; B='254/2'; C=193}
  @{Row=219; A=217; B='337/5'; C=193}
  @{Row=220; A=218; B='393/1'; C=193}
  @{Row=221; A=219; B='393/2'; C=193}
  @{Row=222; A=220; B='393/3'; C=193}
  @{Row=223; A=221; B='465'; C=193}
  @{Row=224; A=222; B='614'; C=193}
  @{Row=225; A=223; B='1303/1'; C=193}
  @{Row=226; A=224; B='1303/2'; C=193}
  @{Row=227; A=225; B='1309'; C=193}
  @{Row=228; A=226; B='1330'; C=193}
  @{Row=229; A=227; B='1334'; C=193}
  @{Row=230; A=228; B='1346'; C=193}
  @{Row=231; A=229; B='1369/1'; C=193}
  @{Row=232; A=230; B='1117/2'; C=193}
  @{Row=233; A=231; B='1230/85'; C=193}
  @{Row=234; A=232; B='1230/86'; C=193}
  @{Row=235; A=233; B='1230/87'; C=193}
  @{Row=236; A=234; B='1230/88'; C=193}
  @{Row=237; A=235; B='1230/100'; C=193}
  @{Row=238; A=236; B='1230/115'; C=193}
  @{Row=239; A=237; B='194/4'; C=193}
)

foreach ($item in $data) {
    $row = $item.Row
    if ($row -ge 224) {
        $ws.Cells.Item($row, 1).Value = $item.A
    }
    $ws.Cells.Item($row, 2).Value = $item.B
    $ws.Cells.Item($row, 3).Value = $item.C
}

# 3) Restore the plain (unstyled) number format on column B now that the
#    values are safely stored as text, so the cells match the rest of the
#    sheet (no extra styling besides text-as-entered).
$fmtSrc = $ws.Range("B2:B6")
$fmtSrc.Copy()
$ws.Range("B158:B239").PasteSpecial(-4122)
